$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.357.88'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '1.881.35'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7120'
$ws.Range("E5").Value = '  +0.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.78'
$ws.Range("E6").Value = '  +0.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.08028'
$ws.Range("E8").Value = '  +3.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3140'
$ws.Range("E9").Value = '  +1.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.12'
$ws.Range("E10").Value = '  -0.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08331'
$ws.Range("E11").Value = '  -1.39%  '
$ws.Range("D12").Value = '1.903.83'
$ws.Range("E12").Value = '  +2.02%  '
$ws.Range("E13").Value = '  +0.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.98'
$ws.Range("E14").Value = '  +4.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7192'
$ws.Range("E15").Value = '  +1.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.366'
$ws.Range("E16").Value = '  +5.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008687'
$ws.Range("D18").Value = '29.381.81'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.69'
$ws.Range("E19").Value = '  +1.11%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.35'
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.147.88'
$ws.Range("E21").Value = '  +1.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.851'
$ws.Range("E23").Value = '  +0.93%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1576'
$ws.Range("E25").Value = '  -2.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.49'
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.095'
$ws.Range("E27").Value = '  +0.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.510'
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.440'
$ws.Range("E30").Value = '  +0.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.363'
$ws.Range("E31").Value = '  +1.47%  '
$ws.Range("E32").Value = '  -6.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05393'
$ws.Range("E33").Value = '  +2.09%  '
$ws.Range("E34").Value = '  +0.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7793'
$ws.Range("E35").Value = '  +4.64%  '
$ws.Range("E36").Value = '  +0.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.688'
$ws.Range("E37").Value = '  -0.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01888'
$ws.Range("E38").Value = '  +1.05%  '
$ws.Range("D39").Value = '1.270.51'
$ws.Range("E39").Value = '  +4.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.746'
$ws.Range("E40").Value = '  +0.89%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.536'
$ws.Range("E41").Value = '  +0.77%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9212'
$ws.Range("E42").Value = '  +3.61%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '113.77'
$ws.Range("E43").Value = '  +4.50%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '74.68'
$ws.Range("E44").Value = '  +2.66%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.002'
$ws.Range("E45").Value = '  +0.14%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '2.063.81'
$ws.Range("E46").Value = '  +2.31%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000128'
$ws.Range("E47").Value = '  +3.56%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.813'
$ws.Range("E48").Value = '  -0.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.5223'
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.575'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4385'
$ws.Range("E51").Value = '  +1.50%  '
